$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Android Application")

# Add the new Android paid course (URL first, then title, to match
# shared-string insertion order used by the original author's save).
$ws.Range("C8").Value = "https://www.udemy.com/course/flutter-bootcamp-with-dart/"
$ws.Range("B8").Value = "The Complete 2021 Flutter Development Bootcamp with Dart"

# Make "Android Application" the active/selected sheet (was "Python").
$ws.Activate()
